# Apply updated "想去人数" (F) and "最低票价" (G) values across the
# "展览" (sheet1), "本地生活" (sheet3) and "全部类型" (sheet4) worksheets,
# matching the refreshed data snapshot described in the commit
# "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" -----------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Cells.Item(4, 6).Value  = 5974
$ws1.Cells.Item(4, 7).Value  = 80
$ws1.Cells.Item(5, 6).Value  = 5974
$ws1.Cells.Item(5, 7).Value  = 80
$ws1.Cells.Item(7, 6).Value  = 3024
$ws1.Cells.Item(7, 7).Value  = 80
$ws1.Cells.Item(8, 6).Value  = 1296
$ws1.Cells.Item(8, 7).Value  = 70
$ws1.Cells.Item(12, 6).Value = 37
$ws1.Cells.Item(13, 6).Value = 325
$ws1.Cells.Item(14, 6).Value = 4452
$ws1.Cells.Item(15, 6).Value = 4452
$ws1.Cells.Item(17, 6).Value = 97
$ws1.Cells.Item(18, 6).Value = 132
$ws1.Cells.Item(21, 6).Value = 81
$ws1.Cells.Item(22, 6).Value = 6901
$ws1.Cells.Item(23, 6).Value = 6901
$ws1.Cells.Item(25, 6).Value = 112
$ws1.Cells.Item(26, 6).Value = 477
$ws1.Cells.Item(29, 6).Value = 1649
$ws1.Cells.Item(32, 6).Value = 6038
$ws1.Cells.Item(38, 6).Value = 6088
$ws1.Cells.Item(39, 6).Value = 16
$ws1.Cells.Item(48, 6).Value = 22
$ws1.Cells.Item(49, 6).Value = 359
$ws1.Cells.Item(52, 6).Value = 1036

# --- Sheet "本地生活" --------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")

$ws3.Cells.Item(2, 6).Value = 1421

# --- Sheet "全部类型" --------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Cells.Item(2, 6).Value  = 1421
$ws4.Cells.Item(4, 6).Value  = 5974
$ws4.Cells.Item(4, 7).Value  = 80
$ws4.Cells.Item(5, 6).Value  = 5974
$ws4.Cells.Item(5, 7).Value  = 80
$ws4.Cells.Item(7, 6).Value  = 3024
$ws4.Cells.Item(7, 7).Value  = 80
$ws4.Cells.Item(8, 6).Value  = 1296
$ws4.Cells.Item(8, 7).Value  = 70
$ws4.Cells.Item(13, 6).Value = 325
$ws4.Cells.Item(14, 6).Value = 4452
$ws4.Cells.Item(15, 6).Value = 4452
$ws4.Cells.Item(17, 6).Value = 97
$ws4.Cells.Item(18, 6).Value = 132
$ws4.Cells.Item(21, 6).Value = 81
$ws4.Cells.Item(22, 6).Value = 6901
$ws4.Cells.Item(23, 6).Value = 6901
$ws4.Cells.Item(25, 6).Value = 112
$ws4.Cells.Item(26, 6).Value = 477
$ws4.Cells.Item(29, 6).Value = 1649
$ws4.Cells.Item(33, 6).Value = 6038
$ws4.Cells.Item(39, 6).Value = 6088
$ws4.Cells.Item(40, 6).Value = 16
$ws4.Cells.Item(48, 6).Value = 359
$ws4.Cells.Item(51, 6).Value = 1036

$wb.Save()
